# mau_bangchamcong_chitiet_chuachot.xlsx
# "Thay doi format muc 7 cho gon hon, in ra theo mau nhu muc 12"
# -> Remove the blank formatting-only row 3 (merged A3:AM3) so the header
#    row (old row 4) becomes row 3, matching the compact layout used for
#    item 12. Also widen a few columns and move the window selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the now-unwanted blank row (old row 3, merged A3:AM3). This
# shifts the header row (old row 4) up to row 3 and fixes dimension /
# mergeCells / row count automatically.
$ws.Rows("3").Delete()

# Widen a few columns to match the new, more compact layout (item 12 style).
$ws.Columns("C").ColumnWidth = 32.85546875
$ws.Columns("D").ColumnWidth = 29.28515625
$ws.Columns("E").ColumnWidth = 25.42578125
$ws.Columns("H").ColumnWidth = 17.140625

# Move the active selection as saved in the edited workbook.
$ws.Range("D12").Select() | Out-Null
